$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Generator Data" -- add a second generator column (C) and
# update the existing generator-1 values.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Generator Data")

# New header for the second generator, copying the formatting (bold,
# border, centered) that "Generator 1" already uses in B1.
$ws1.Range("C1").Value = "Generator 2"
$ws1.Range("B1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)

# Updated Generator 1 numbers (column B).
$ws1.Range("B2").Value = 36090.1058986
$ws1.Range("B3").Value = 7218.02117972
$ws1.Range("B4").Value = 721.802117972
$ws1.Range("B5").Value = 222682.17895

# New Generator 2 numbers (column C).
$ws1.Range("C2").Value = 0.7734556997549999
$ws1.Range("C3").Value = 0.9281468397059999
$ws1.Range("C4").Value = 0.4176660778676999
$ws1.Range("C5").Value = 4.40044706213

# ---------------------------------------------------------------------
# Sheet 2: "Yearly Fuel Costs" -- extend the scenario from 4 years to 20
# years (rows 2-5 existing, rows 6-21 new), and tweak the existing
# values slightly.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

$ws2.Range("B2").Value = 35576.02469177476
$ws2.Range("B3").Value = 35576.09666751218
$ws2.Range("B4").Value = 35576.17640765515
$ws2.Range("B5").Value = 35576.26828023807

$yearlyValues = @(
    35576.37178682209,
    35576.4867560025,
    35576.61588381077,
    35576.76149660935,
    35576.92414064689,
    35577.10183779884,
    35577.29656010856,
    35577.51278889657,
    35577.75446436717,
    35578.0257721881,
    35578.33069279457,
    35578.67413798285,
    35579.06005444584,
    35579.49298712097,
    35579.98074800181,
    35580.5291117245
)

# Use row 2 (A2) as the formatting template for the new label cells, the
# same way the header was extended on sheet 1 above.
$ws2.Range("A2").Copy()

for ($i = 0; $i -lt $yearlyValues.Length; $i++) {
    $row = 6 + $i
    $year = 5 + $i
    $ws2.Range("A$row").Value = "Total Fuel Cost at y = $year"
    $ws2.Range("A$row").PasteSpecial(-4122)
    $ws2.Range("B$row").Value = $yearlyValues[$i]
}
